# MOD_5_var_AV-MCPS.xlsx — add a new "d=6" row between the existing
# "d=5" (row 6) and "d=7" (row 7) rows, pushing "d=7" and "d=10" down
# by one row (so the old row 7 -> row 8, old row 8 -> row 9), and fill
# the new row with the corrected Diebold-Mariano values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Carry the bold/centered/bordered label formatting (style of the
#        existing "d=..." cells in column A) onto the two row-label cells
#        that will end up holding text after the shift: the brand new
#        row 7 ("d=6") and the row that becomes row 9 ("d=10"), which
#        starts out blank/unformatted. Formats only, so no values move yet.
$ws.Range("A6").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)

# --- 2. Shift the old row 8 ("d=10") down into row 9.
$ws.Range("A9").Value = $ws.Range("A8").Value()
$ws.Range("B9").Value = $ws.Range("B8").Value()
$ws.Range("C9").Value = $ws.Range("C8").Value()
$ws.Range("D9").Value = $ws.Range("D8").Value()
$ws.Range("E9").Value = $ws.Range("E8").Value()

# --- 3. Shift the old row 7 ("d=7") down into row 8.
$ws.Range("A8").Value = $ws.Range("A7").Value()
$ws.Range("B8").Value = $ws.Range("B7").Value()
$ws.Range("C8").Value = $ws.Range("C7").Value()
$ws.Range("D8").Value = $ws.Range("D7").Value()
$ws.Range("E8").Value = $ws.Range("E7").Value()

# --- 4. Fill the now-free row 7 with the new "d=6" data.
$ws.Range("A7").Value = "d=6"
$ws.Range("B7").Value = 97.93361859143396
$ws.Range("C7").Value = 98.00158872726553
$ws.Range("D7").Value = 98.04165037344067
$ws.Range("E7").Value = 97.96989947472225

Write-Output "Inserted d=6 row; dimension now A1:E9"
